$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A70").Value = "NEAR-USD"
$ws.Range("A71").Value = "GRT-USD"
